$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current last data row (row 8), pushing the
# existing row 8 down to row 10 (with its formatting/style intact).
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(8).Insert()

# New row 8: updated weekly record (same lot, new date, 12kg granel presentation)
$ws.Cells.Item(8, 1).Value2 = 11
$ws.Cells.Item(8, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(8, 3).Value2 = "Bíobío"
$ws.Cells.Item(8, 4).Value2 = 45077
$ws.Cells.Item(8, 4).NumberFormat = $ws.Cells.Item(7, 4).NumberFormat
$ws.Cells.Item(8, 5).Value2 = 8
$ws.Cells.Item(8, 6).Value2 = "Fruta"
$ws.Cells.Item(8, 7).Value2 = 100107
$ws.Cells.Item(8, 8).Value2 = "Otros"
$ws.Cells.Item(8, 9).Value2 = 100107001
$ws.Cells.Item(8, 10).Value2 = "Caqui"
$ws.Cells.Item(8, 11).Value2 = "Mankaki"
$ws.Cells.Item(8, 12).Value2 = "Primera"
$ws.Cells.Item(8, 13).Value2 = 140
$ws.Cells.Item(8, 14).Value2 = 12000
$ws.Cells.Item(8, 15).Value2 = 14000
$ws.Cells.Item(8, 16).Value2 = 12857
$ws.Cells.Item(8, 17).Value2 = "$/caja 12 kilos granel"
$ws.Cells.Item(8, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(8, 19).Value2 = 12857
$ws.Cells.Item(8, 20).Value2 = 1

# New row 9: new weekly record, second quality
$ws.Cells.Item(9, 1).Value2 = 11
$ws.Cells.Item(9, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(9, 3).Value2 = "Bíobío"
$ws.Cells.Item(9, 4).Value2 = 45077
$ws.Cells.Item(9, 4).NumberFormat = $ws.Cells.Item(7, 4).NumberFormat
$ws.Cells.Item(9, 5).Value2 = 8
$ws.Cells.Item(9, 6).Value2 = "Fruta"
$ws.Cells.Item(9, 7).Value2 = 100107
$ws.Cells.Item(9, 8).Value2 = "Otros"
$ws.Cells.Item(9, 9).Value2 = 100107001
$ws.Cells.Item(9, 10).Value2 = "Caqui"
$ws.Cells.Item(9, 11).Value2 = "Mankaki"
$ws.Cells.Item(9, 12).Value2 = "Segunda"
$ws.Cells.Item(9, 13).Value2 = 80
$ws.Cells.Item(9, 14).Value2 = 11000
$ws.Cells.Item(9, 15).Value2 = 11000
$ws.Cells.Item(9, 16).Value2 = 11000
$ws.Cells.Item(9, 17).Value2 = "$/caja 12 kilos granel"
$ws.Cells.Item(9, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(9, 19).Value2 = 11000
$ws.Cells.Item(9, 20).Value2 = 1
